$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.311.70'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.538.03'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'608.24"
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = "'143.97"
$ws.Range("E6").Value = '  -2.67%  '
$ws.Range("D7").Value = '3.536.73'
$ws.Range("E7").Value = '  +0.67%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("E11").Value = '  -4.33%  '
$ws.Range("E12").Value = '  -2.78%  '
$ws.Range("D13").Value = '4.137.89'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").Value = "'0.0000208"
$ws.Range("E14").Value = '  -4.85%  '
$ws.Range("D15").Value = "'30.27"
$ws.Range("E15").Value = '  -5.29%  '
$ws.Range("D16").Value = '3.533.12'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '66.403.63'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").Value = "'10.95"
$ws.Range("E19").Value = '  +1.81%  '
$ws.Range("E20").Value = '  -3.87%  '
$ws.Range("D21").Value = "'14.96"
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("D22").Value = "'425.89"
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("D23").Value = "'0.602"
$ws.Range("E23").Value = '  -1.29%  '
$ws.Range("D24").Value = "'78.69"
$ws.Range("E24").Value = '  -1.10%  '
$ws.Range("D25").Value = '3.679.85'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("D28").Value = "'8.11"
$ws.Range("E28").Value = '  -2.22%  '
$ws.Range("D29").Value = "'9.21"
$ws.Range("E29").Value = '  -5.96%  '
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E32").Value = '  -8.18%  '
$ws.Range("E33").Value = '  -4.47%  '
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("D35").Value = '3.528.18'
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  -2.94%  '
$ws.Range("D38").Value = "'5.65"
$ws.Range("E38").Value = '  -5.46%  '
$ws.Range("D39").Value = "'7.83"
$ws.Range("E39").Value = '  -2.77%  '
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").Value = "'172.54"
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("E42").Value = '  -4.32%  '
$ws.Range("E43").Value = '  -4.26%  '
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("E45").Value = '  -7.97%  '
$ws.Range("D46").Value = "'45.52"
$ws.Range("E46").Value = '  -1.42%  '
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("D48").Value = "'26.10"
$ws.Range("E48").Value = '  -6.93%  '
$ws.Range("D49").Value = "'2.41"
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("E50").Value = '  -4.40%  '
$ws.Range("D51").Value = "'0.947"
$ws.Range("E51").Value = '  -4.58%  '
